$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) values are written as text, matching the
# source data (avoids Excel auto-coercing numeric-looking strings
# like "238.82" into floating point numbers, which would corrupt
# values such as "1.900.50" or drop trailing zeros like "7.380").

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.485.63"
$ws.Range("E2").Value = "  +0.27%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.900.50"
$ws.Range("E3").Value = "  -0.29%  "

$ws.Range("E4").Value = "  +0.28%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "238.82"
$ws.Range("E5").Value = "  +0.09%  "

$ws.Range("E6").Value = "  +0.23%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4723"
$ws.Range("E7").Value = "  -0.97%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2847"
$ws.Range("E8").Value = "  -0.52%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06638"
$ws.Range("E9").Value = "  -0.90%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.67"
$ws.Range("E10").Value = "  +4.49%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "99.36"
$ws.Range("E11").Value = "  -3.99%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07798"
$ws.Range("E12").Value = "  +0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.908.16"
$ws.Range("E13").Value = "  +0.18%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.170"
$ws.Range("E14").Value = "  -0.87%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6759"
$ws.Range("E15").Value = "  -0.18%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "284.31"
$ws.Range("E16").Value = "  +8.22%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.479.13"
$ws.Range("E17").Value = "  +0.21%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.34%  "

$ws.Range("B19").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C19").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.161.80"
$ws.Range("E19").Value = "  +0.41%  "

$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007456"
$ws.Range("E20").Value = "  -0.77%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.69"
$ws.Range("E21").Value = "  -0.21%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.394"
$ws.Range("E22").Value = "  -1.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  +0.60%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.266"
$ws.Range("E24").Value = "  -0.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.348"
$ws.Range("E25").Value = "  -1.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "167.27"
$ws.Range("E26").Value = "  +1.73%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.20"
$ws.Range("E27").Value = "  +0.56%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.007"
$ws.Range("E28").Value = "  -2.64%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.392"
$ws.Range("E29").Value = "  +1.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09899"
$ws.Range("E30").Value = "  -2.42%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.516"
$ws.Range("E31").Value = "  -2.98%  "

$ws.Range("E32").Value = "  -0.30%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.249"
$ws.Range("E33").Value = "  -0.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04737"
$ws.Range("E34").Value = "  -1.01%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7179"
$ws.Range("E35").Value = "  -2.13%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.106"
$ws.Range("E36").Value = "  -0.89%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.718"
$ws.Range("E37").Value = "  +0.38%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01888"
$ws.Range("E38").Value = "  -1.95%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.725"
$ws.Range("E39").Value = "  +7.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.555"
$ws.Range("E40").Value = "  -1.46%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.44"
$ws.Range("E41").Value = "  -1.83%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.980"
$ws.Range("E42").Value = "  -0.86%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.8685"
$ws.Range("E43").Value = "  +0.59%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "104.90"
$ws.Range("E44").Value = "  -1.56%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4255"
$ws.Range("E45").Value = "  -0.30%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.003"
$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "987.77"
$ws.Range("E47").Value = "  -3.07%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.380"
$ws.Range("E48").Value = "  -1.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.219"
$ws.Range("E49").Value = "  +4.16%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.1180"
$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.39"
$ws.Range("E51").Value = "  -2.17%  "
